$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $jCell = $ws.Cells.Item($r, 10)  # Column J (level_1)
    if ($jCell.Value2 -eq "JV") {
        $jCell.Value = "Junior Varsity"
    }
    $kCell = $ws.Cells.Item($r, 11)  # Column K (level_2)
    if ($kCell.Value2 -eq "JV") {
        $kCell.Value = "Junior Varsity"
    }
}
